$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin price/link/volume cells are stored as plain text in this sheet (e.g. "1.003",
# "29.118.40"). Assigning a numeric-looking string straight to Range.Value lets Excel
# auto-convert it to a number, which would change the cell's underlying type/format
# and the diff only ever touches text. Set-TextCell keeps the text that way: numbers
# made only of digits and a single decimal point are written with a leading
# apostrophe so Excel stores/keeps them as literal text (matching the original
# inline-string cells), everything else (names, URLs, the "  +x.xx%  " volumes,
# and multi-dot price strings such as "29.206.27") is assigned as-is.
function Set-TextCell($addr, $text) {
    $value = $text
    if ($text -match '^[0-9]+\.[0-9]+$') {
        $value = "'" + $text
    }
    $ws.Range($addr).Value = $value
}

Set-TextCell 'D2' '29.206.27'
Set-TextCell 'E2' '  +0.60%  '

Set-TextCell 'D3' '1.897.01'
Set-TextCell 'E3' '  +0.36%  '

Set-TextCell 'D4' '1.001'
Set-TextCell 'E4' '  +0.17%  '

Set-TextCell 'D5' '323.71'
Set-TextCell 'E5' '  -1.80%  '

Set-TextCell 'D6' '1.001'
Set-TextCell 'E6' '  +0.26%  '

Set-TextCell 'D7' '0.4706'
Set-TextCell 'E7' '  +2.74%  '

Set-TextCell 'D8' '0.4020'
Set-TextCell 'E8' '  -2.19%  '

Set-TextCell 'D9' '47.53'
Set-TextCell 'E9' '  -0.35%  '

Set-TextCell 'D10' '0.07993'
Set-TextCell 'E10' '  +0.33%  '

Set-TextCell 'D11' '0.9931'
Set-TextCell 'E11' '  -0.45%  '

Set-TextCell 'D12' '22.53'
Set-TextCell 'E12' '  +3.59%  '

Set-TextCell 'D13' '1.902.27'
Set-TextCell 'E13' '  +0.88%  '

Set-TextCell 'D14' '5.855'
Set-TextCell 'E14' '  -1.11%  '

Set-TextCell 'D15' '7.040'
Set-TextCell 'E15' '  -0.69%  '

Set-TextCell 'B16' 'Litecoin'
Set-TextCell 'C16' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D16' '88.83'
Set-TextCell 'E16' '  +0.20%  '

Set-TextCell 'B17' 'BinanceUSD'
Set-TextCell 'C17' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextCell 'D17' '1.002'
Set-TextCell 'E17' '  +0.18%  '

Set-TextCell 'D18' '0.06613'
Set-TextCell 'E18' '  +0.79%  '

Set-TextCell 'D19' '0.00001025'
Set-TextCell 'E19' '  -0.05%  '

Set-TextCell 'D20' '17.47'
Set-TextCell 'E20' '  +0.17%  '

Set-TextCell 'E21' '  +0.42%  '

Set-TextCell 'D22' '29.223.91'
Set-TextCell 'E22' '  +0.64%  '

Set-TextCell 'D23' '5.494'
Set-TextCell 'E23' '  +1.08%  '

Set-TextCell 'D24' '11.54'
Set-TextCell 'E24' '  +0.88%  '

Set-TextCell 'D25' '2.198'
Set-TextCell 'E25' '  -0.09%  '

Set-TextCell 'D26' '2.111.23'
Set-TextCell 'E26' '  -0.17%  '

Set-TextCell 'D27' '153.97'
Set-TextCell 'E27' '  -1.42%  '

Set-TextCell 'D28' '19.63'
Set-TextCell 'E28' '  +0.20%  '

Set-TextCell 'D29' '5.993'
Set-TextCell 'E29' '  +8.90%  '

Set-TextCell 'D30' '2.087'
Set-TextCell 'E30' '  -0.11%  '

Set-TextCell 'D31' '117.23'
Set-TextCell 'E31' '  -0.24%  '

Set-TextCell 'D32' '1.062'
Set-TextCell 'E32' '  +2.53%  '

Set-TextCell 'D33' '0.09446'
Set-TextCell 'E33' '  +1.33%  '

Set-TextCell 'B34' 'ARBITRUM'
Set-TextCell 'C34' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell 'D34' '1.398'
Set-TextCell 'E34' '  -0.77%  '

Set-TextCell 'B35' 'HuobiToken'
Set-TextCell 'C35' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D35' '3.541'
Set-TextCell 'E35' '  +0.41%  '

Set-TextCell 'D36' '5.341'
Set-TextCell 'E36' '  +0.79%  '

Set-TextCell 'D37' '0.06068'
Set-TextCell 'E37' '  +0.19%  '

Set-TextCell 'D38' '0.02244'
Set-TextCell 'E38' '  +0.72%  '

Set-TextCell 'E39' '  -0.05%  '

Set-TextCell 'D40' '8.062'
Set-TextCell 'E40' '  -3.75%  '

Set-TextCell 'D41' '0.5804'
Set-TextCell 'E41' '  +0.23%  '

Set-TextCell 'D42' '0.1824'
Set-TextCell 'E42' '  -0.11%  '

Set-TextCell 'D43' '2.473'
Set-TextCell 'E43' '  +8.01%  '

Set-TextCell 'D44' '10.04'
Set-TextCell 'E44' '  -0.66%  '

Set-TextCell 'B45' 'WEMIXToken'
Set-TextCell 'C45' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell 'D45' '1.268'
Set-TextCell 'E45' '  +0.78%  '

Set-TextCell 'B46' 'Cronos'
Set-TextCell 'C46' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextCell 'D46' '0.07703'
Set-TextCell 'E46' '  +2.60%  '

Set-TextCell 'D47' '12.16'
Set-TextCell 'E47' '  +1.10%  '

Set-TextCell 'D48' '0.5468'
Set-TextCell 'E48' '  +0.05%  '

Set-TextCell 'D49' '1.899'
Set-TextCell 'E49' '  -0.20%  '

Set-TextCell 'D50' '113.27'
Set-TextCell 'E50' '  +1.97%  '

Set-TextCell 'D51' '43.68'
Set-TextCell 'E51' '  -2.00%  '
